$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Python")
$ws2 = $wb.Worksheets.Item("cert_data")

# --- cert_data sheet: add a new "без модуля" entry under the module column (B) ---
$ws2.Range("B4").Value = "без модуля"
# Match the formatting used by the other cells in that column/area (style s="8").
$ws2.Range("B3").Copy()
$ws2.Range("B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Python sheet: reassign module/certificate values on rows 2-4 ---
$ws1.Range("E2").Value = "без модуля"
$ws1.Range("E3").Value = "базовый модуль"
$ws1.Range("G3").Value = "сертификат с отличием"
$ws1.Range("E4").Value = "углубленный модуль"

# --- Selections recorded on each sheet's active view ---
[void]$ws2.Range("B4").Select()
[void]$ws1.Range("G3").Select()
[void]$ws1.Activate()
